$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added for "Jengibre" at Vega Modelo de Temuco.
# Insert a new row at row 60, pushing the existing rows 60-102 down to 61-103.
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new record's data.
$ws.Range("A60").Value = 10
$ws.Range("B60").Value = "Vega Modelo de Temuco"
$ws.Range("C60").Value = "La Araucanía"
$ws.Range("D60").Value = 44438
$ws.Range("E60").Value = 9
$ws.Range("F60").Value = 100114007
$ws.Range("G60").Value = "Jengibre"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 20
$ws.Range("K60").Value = 30000
$ws.Range("L60").Value = 30000
$ws.Range("M60").Value = 30000
$ws.Range("N60").Value = "$/caja 13 kilos"
$ws.Range("O60").Value = "Perú"
$ws.Range("P60").Value = 2308
$ws.Range("Q60").Value = 13
$ws.Range("R60").Value = "Hortaliza"

# Match the date-column number format (yyyy-mm-dd hh:mm:ss) used by the rest
# of column D by copying the style from the row right below.
$ws.Range("D61").Copy()
$ws.Range("D60").PasteSpecial(-4122)
